$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("ECs", "Lpl", "Vldlr", "ECs", 3, 1, 71.329076, 213.987228, 0.3307464087015077, 0.3307464087015077, 2, 0.6666666666666666, 0.3374003333333333, 1.012201, 0.01738364872808817, 0.01738364872808818, 24.06645401875867, 216.598086168828, 0.005749579386943695, 0.005749579386943698)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

$row3 = @("ECs", "Lpl", "Vldlr", "FAPs", 3, 1, 71.329076, 213.987228, 0.3307464087015077, 0.3307464087015077, 3, 1, 16.71131166666667, 50.133935, 0.8610055862391021, 0.8610055862391023, 1192.002419931353, 10728.02177938218, 0.2847745055205192, 0.2847745055205194)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

$row4 = @("FAPs", "Lpl", "Vldlr", "sCs", 3, 1, 71.329076, 213.987228, 0.3307464087015077, 0.3307464087015077, 3, 1, 2.360351, 7.081053000000001, 0.1216107650328097, 0.1216107650328097, 168.361655865676, 1515.254902791084, 0.04022232379404468, 0.0402223237940447)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}

$row5 = @("FAPs", "Lpl", "Vldlr", "ECs", 3, 1, 129.5615336666667, 388.684601, 0.6007649947142101, 0.6007649947142101, 2, 0.6666666666666666, 0.3374003333333333, 1.012201, 0.01738364872808817, 0.01738364872808818, 43.71410464631123, 393.426941816801, 0.01044348763624358, 0.01044348763624358)
for ($i = 0; $i -lt $row5.Length; $i++) {
    $ws.Cells.Item(5, $i + 1).Value = $row5[$i]
}

$row6 = @("sCs", "Lpl", "Vldlr", "FAPs", 3, 1, 129.5615336666667, 388.684601, 0.6007649947142101, 0.6007649947142101, 3, 1, 16.71131166666667, 50.133935, 0.8610055862391021, 0.8610055862391023, 2165.143169114993, 19486.28852203494, 0.5172620164658395, 0.5172620164658396)
for ($i = 0; $i -lt $row6.Length; $i++) {
    $ws.Cells.Item(6, $i + 1).Value = $row6[$i]
}

$row7 = @("sCs", "Lpl", "Vldlr", "sCs", 3, 1, 129.5615336666667, 388.684601, 0.6007649947142101, 0.6007649947142101, 3, 1, 2.360351, 7.081053000000001, 0.1216107650328097, 0.1216107650328097, 305.8106955516504, 2752.296259964854, 0.07305949061212694, 0.07305949061212696)
for ($i = 0; $i -lt $row7.Length; $i++) {
    $ws.Cells.Item(7, $i + 1).Value = $row7[$i]
}

$row8 = @("sCs", "Lpl", "Vldlr", "ECs", 3, 1, 14.770314, 44.310942, 0.06848859658428214, 0.06848859658428215, 2, 0.6666666666666666, 0.3374003333333333, 1.012201, 0.01738364872808817, 0.01738364872808818, 4.983508867038, 44.851579803342, 0.0011905817049009, 0.001190581704900901)
for ($i = 0; $i -lt $row8.Length; $i++) {
    $ws.Cells.Item(8, $i + 1).Value = $row8[$i]
}

$row9 = @("sCs", "Lpl", "Vldlr", "FAPs", 3, 1, 14.770314, 44.310942, 0.06848859658428214, 0.06848859658428215, 3, 1, 16.71131166666667, 50.133935, 0.8610055862391021, 0.8610055862391023, 246.83132066853, 2221.48188601677, 0.0589690642527432, 0.05896906425274323)
for ($i = 0; $i -lt $row9.Length; $i++) {
    $ws.Cells.Item(9, $i + 1).Value = $row9[$i]
}

$row10 = @("sCs", "Lpl", "Vldlr", "sCs", 3, 1, 14.770314, 44.310942, 0.06848859658428214, 0.06848859658428215, 3, 1, 2.360351, 7.081053000000001, 0.1216107650328097, 0.1216107650328097, 34.863125420214, 313.7681287819261, 0.008328950626638024, 0.00832895062663803)
for ($i = 0; $i -lt $row10.Length; $i++) {
    $ws.Cells.Item(10, $i + 1).Value = $row10[$i]
}

